# Add a new worksheet "ODI Batting Extra" after the existing sheets,
# populated with extra per-match batting stats, mirroring the header
# style used on the other sheets.

$wb = $excel.ActiveWorkbook

# Source sheet whose header formatting (bold, border, centered) we reuse
# for the new sheet's header row, so we don't create a brand-new style.
$wsSource = $wb.Worksheets.Item(2)

# Insert the new sheet at the end of the workbook (after the last sheet).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "ODI Batting Extra"

# Header row.
$headers = @("MATCH_CODE","BATTING_POSITION","NUM_4","NUM_6","PERCENT_RUNS_OF_TOTAL","MAN_OF_MATCH")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Copy the header formatting (bold font, borders, centered alignment) from
# an existing header cell so the new header row matches the workbook's
# existing look instead of inventing a new style.
$wsSource.Range("A1").Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)

# Data rows. Values that look numeric (match codes, the "0" placeholders,
# the percentage string) need to stay text, so they are entered with a
# leading apostrophe to force text entry; BATTING_POSITION on row 3 is a
# genuine number.
$ws.Range("A2").Value = "'4727"
$ws.Range("B2").Value = "'"
$ws.Range("C2").Value = "'"
$ws.Range("D2").Value = "'"
$ws.Range("E2").Value = "'"
$ws.Range("F2").Value = "NO"

$ws.Range("A3").Value = "'4731"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = "'0"
$ws.Range("D3").Value = "'0"
$ws.Range("E3").Value = "'1.14%"
$ws.Range("F3").Value = "NO"

# Drop the "quote prefix" styling that forcing text entry leaves behind on
# the data cells, so they fall back to the default (unstyled) look.
$ws.Range("A2:F3").Style = "Normal"
